$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 129
$ws.Range("I2").Value = 129
$ws.Range("K2").Value = 129
$ws.Range("M2").Value = -16
$ws.Range("H5").Value = 160.6
$ws.Range("I5").Value = 76.5
$ws.Range("J5").Value = 497
$ws.Range("K5").Value = 76.5
$ws.Range("L5").Value = 497
$ws.Range("M5").Value = 38.5
$ws.Range("N5").Value = -727
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = $null
$ws.Range("N43").Value = -1138
$ws.Range("H70").Value = 1854.4445
$ws.Range("J70").Value = 1855.7142
$ws.Range("L70").Value = 5567.142599999999
$ws.Range("N70").Value = -6107.142599999999
$ws.Range("H73").Value = 1854.4445
$ws.Range("J73").Value = 1855.7142
$ws.Range("L73").Value = 5567.142599999999
$ws.Range("N73").Value = -7439.142599999999
$ws.Range("H76").Value = 3328.3333
$ws.Range("I76").Value = 2100
$ws.Range("K76").Value = 2100
$ws.Range("M76").Value = -1785
$ws.Range("H79").Value = 3328.3333
$ws.Range("I79").Value = 2100
$ws.Range("K79").Value = 2100
$ws.Range("M79").Value = -1008
$ws.Range("H80").Value = 2587.1875
$ws.Range("I80").Value = 2683.75
$ws.Range("J80").Value = 2490.625
$ws.Range("K80").Value = 8051.25
$ws.Range("L80").Value = 7471.875
$ws.Range("M80").Value = -7053.25
$ws.Range("N80").Value = -9467.875
$ws.Range("H83").Value = 2587.1875
$ws.Range("I83").Value = 2683.75
$ws.Range("J83").Value = 2490.625
$ws.Range("K83").Value = 24153.75
$ws.Range("L83").Value = 22415.625
$ws.Range("M83").Value = -19161.75
$ws.Range("N83").Value = -32399.625
$ws.Range("H93").Value = 42500
$ws.Range("J93").Value = 42500
$ws.Range("L93").Value = 42500
$ws.Range("N93").Value = -47492
$ws.Range("H94").Value = 1524
$ws.Range("I94").Value = 1098
$ws.Range("J94").Value = 1950
$ws.Range("K94").Value = 1098
$ws.Range("L94").Value = 1950
$ws.Range("M94").Value = -647
$ws.Range("N94").Value = -2852
$ws.Range("H96").Value = 399.53845
$ws.Range("I96").Value = 286.55554
$ws.Range("K96").Value = 859.66662
$ws.Range("M96").Value = 513.33338
$ws.Range("H99").Value = 1049.5
$ws.Range("I99").Value = 1049.5
$ws.Range("K99").Value = 3148.5
$ws.Range("M99").Value = -1650.5
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 10000
$ws.Range("K113").Value = 10000
$ws.Range("M113").Value = -6746
$ws.Range("H116").Value = 4499.2
$ws.Range("I116").Value = 4499.25
$ws.Range("K116").Value = 4499.25
$ws.Range("M116").Value = -1057.25
$ws.Range("H132").Value = 1294.3529
$ws.Range("I132").Value = 1294.3529
$ws.Range("K132").Value = 3883.0587
$ws.Range("M132").Value = -1353.0587
$ws.Range("H138").Value = 3273.2188
$ws.Range("J138").Value = 3843.5417
$ws.Range("L138").Value = 11530.6251
$ws.Range("N138").Value = -21810.6251
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12059.5625
$ws.Range("I32").Value = 12059.5625
$ws.Range("K32").Value = 12059.5625
$ws.Range("M32").Value = -11772.5625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1246
$ws.Range("I99").Value = 1246
$ws.Range("K99").Value = 1246
$ws.Range("M99").Value = 252
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2577
$ws.Range("I31").Value = 2662.6667
$ws.Range("K31").Value = 2662.6667
$ws.Range("M31").Value = -2367.6667
$ws.Range("H34").Value = 2577
$ws.Range("I34").Value = 2662.6667
$ws.Range("K34").Value = 2662.6667
$ws.Range("M34").Value = -2460.6667
$ws.Range("H93").Value = 35000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 35000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 35000
$ws.Range("M93").Value = $null
$ws.Range("N93").Value = -38744
$ws.Range("H107").Value = 2135.9285
$ws.Range("I107").Value = 2373.5454
$ws.Range("J107").Value = 1264.6666
$ws.Range("K107").Value = 2373.5454
$ws.Range("L107").Value = 1264.6666
$ws.Range("M107").Value = -453.5454
$ws.Range("N107").Value = -5104.6666
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 76.28570999999999
$ws.Range("J2").Value = 20
$ws.Range("L2").Value = 120
$ws.Range("N2").Value = -346
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("H6").Value = 176.125
$ws.Range("I6").Value = 176.125
$ws.Range("K6").Value = 528.375
$ws.Range("M6").Value = -415.375
$ws.Range("H44").Value = 503
$ws.Range("I44").Value = 503
$ws.Range("K44").Value = 1509
$ws.Range("M44").Value = -1111
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 286.6
$ws.Range("I2").Value = 176.42857
$ws.Range("K2").Value = 176.42857
$ws.Range("M2").Value = -63.42857000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 4002256
$ws.Range("I19").Value = 5002800
$ws.Range("J19").Value = 81
$ws.Range("K19").Value = 5002800
$ws.Range("L19").Value = 81
$ws.Range("M19").Value = -5002630
$ws.Range("N19").Value = -421
$ws.Range("H22").Value = 4520
$ws.Range("I22").Value = 9833.333000000001
$ws.Range("J22").Value = 2242.8572
$ws.Range("K22").Value = 9833.333000000001
$ws.Range("L22").Value = 2242.8572
$ws.Range("M22").Value = -9538.333000000001
$ws.Range("N22").Value = -2832.8572
$ws.Range("H27").Value = 4520
$ws.Range("I27").Value = 9833.333000000001
$ws.Range("J27").Value = 2242.8572
$ws.Range("K27").Value = 9833.333000000001
$ws.Range("L27").Value = 2242.8572
$ws.Range("M27").Value = -9726.333000000001
$ws.Range("N27").Value = -2456.8572
$ws.Range("H61").Value = 1024.25
$ws.Range("I61").Value = 699
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 699
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -497
$ws.Range("N61").Value = -2404
$ws.Range("H113").Value = 1024.25
$ws.Range("I113").Value = 699
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 699
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1471
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 6633
$ws.Range("I122").Value = 6633
$ws.Range("K122").Value = 19899
$ws.Range("M122").Value = -17449
$ws.Range("H132").Value = 4963.8945
$ws.Range("I132").Value = 3815.4443
$ws.Range("K132").Value = 11446.3329
$ws.Range("M132").Value = -8916.332900000001
$ws.Range("H136").Value = 3685.875
$ws.Range("I136").Value = 3581.1667
$ws.Range("K136").Value = 10743.5001
$ws.Range("M136").Value = -8193.500100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1514.5
$ws.Range("I100").Value = 603.4
$ws.Range("J100").Value = 3033
$ws.Range("K100").Value = 1206.8
$ws.Range("L100").Value = 6066
$ws.Range("M100").Value = -665.8
$ws.Range("N100").Value = -7148
$ws.Range("H107").Value = 1108
$ws.Range("I107").Value = 756.7143
$ws.Range("J107").Value = 1599.8
$ws.Range("K107").Value = 2270.1429
$ws.Range("L107").Value = 4799.4
$ws.Range("M107").Value = -350.1428999999998
$ws.Range("N107").Value = -8639.4
$ws.Range("H113").Value = 566.5454999999999
$ws.Range("I113").Value = 552.1667
$ws.Range("K113").Value = 1656.5001
$ws.Range("M113").Value = 513.4999
